$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension implicitly handled by writing to A1:T7 range; fill rows 2-7

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Fgf16"
$ws.Range("C2").Value = "Fgfr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.678104
$ws.Range("H2").Value = 5.034312
$ws.Range("I2").Value = 0.551436927751233
$ws.Range("J2").Value = 0.551436927751233
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.5555316666666666
$ws.Range("N2").Value = 1.666595
$ws.Range("O2").Value = 0.01938483203642842
$ws.Range("P2").Value = 0.01938483203642843
$ws.Range("Q2").Value = 0.93223991196
$ws.Range("R2").Value = 8.39015920764
$ws.Range("S2").Value = 0.01068951222314177
$ws.Range("T2").Value = 0.01068951222314177

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Fgf16"
$ws.Range("C3").Value = "Fgfr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.678104
$ws.Range("H3").Value = 5.034312
$ws.Range("I3").Value = 0.551436927751233
$ws.Range("J3").Value = 0.551436927751233
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1817723333333333
$ws.Range("N3").Value = 0.545317
$ws.Range("O3").Value = 0.006342799811357313
$ws.Range("P3").Value = 0.006342799811357313
$ws.Range("Q3").Value = 0.305032879656
$ws.Range("R3").Value = 2.745295916904
$ws.Range("S3").Value = 0.003497654041315976
$ws.Range("T3").Value = 0.003497654041315977

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fgf16"
$ws.Range("C4").Value = "Fgfr4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.678104
$ws.Range("H4").Value = 5.034312
$ws.Range("I4").Value = 0.551436927751233
$ws.Range("J4").Value = 0.551436927751233
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 27.92075533333333
$ws.Range("N4").Value = 83.762266
$ws.Range("O4").Value = 0.9742723681522142
$ws.Range("P4").Value = 0.9742723681522143
$ws.Range("Q4").Value = 46.853931207888
$ws.Range("R4").Value = 421.685380870992
$ws.Range("S4").Value = 0.5372497614867752
$ws.Range("T4").Value = 0.5372497614867753

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Fgf16"
$ws.Range("C5").Value = "Fgfr4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.365043666666667
$ws.Range("H5").Value = 4.095131
$ws.Range("I5").Value = 0.448563072248767
$ws.Range("J5").Value = 0.448563072248767
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.5555316666666666
$ws.Range("N5").Value = 1.666595
$ws.Range("O5").Value = 0.01938483203642842
$ws.Range("P5").Value = 0.01938483203642843
$ws.Range("Q5").Value = 0.7583249832161112
$ws.Range("R5").Value = 6.824924848945001
$ws.Range("S5").Value = 0.008695319813286655
$ws.Range("T5").Value = 0.008695319813286658

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Fgf16"
$ws.Range("C6").Value = "Fgfr4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.365043666666667
$ws.Range("H6").Value = 4.095131
$ws.Range("I6").Value = 0.448563072248767
$ws.Range("J6").Value = 0.448563072248767
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1817723333333333
$ws.Range("N6").Value = 0.545317
$ws.Range("O6").Value = 0.006342799811357313
$ws.Range("P6").Value = 0.006342799811357313
$ws.Range("Q6").Value = 0.2481271723918889
$ws.Range("R6").Value = 2.233144551527
$ws.Range("S6").Value = 0.002845145770041336
$ws.Range("T6").Value = 0.002845145770041337

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Fgf16"
$ws.Range("C7").Value = "Fgfr4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.365043666666667
$ws.Range("H7").Value = 4.095131
$ws.Range("I7").Value = 0.448563072248767
$ws.Range("J7").Value = 0.448563072248767
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 27.92075533333333
$ws.Range("N7").Value = 83.762266
$ws.Range("O7").Value = 0.9742723681522142
$ws.Range("P7").Value = 0.9742723681522143
$ws.Range("Q7").Value = 38.11305023631623
$ws.Range("R7").Value = 343.017452126846
$ws.Range("S7").Value = 0.4370226066654389
$ws.Range("T7").Value = 0.437022606665439

